# Apply the "BNEF" source-citation swap to the PDiBCpDoC workbook.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("PDiBCpDoC")

# --- "About" sheet: replace the MIT/RSC citation with the BNEF citation ---
$ws1.Range("B3").Value = "BNEF"
$ws1.Range("B4").Value = 2019
$ws1.Range("B5").Value = "A Behind the Scenes Take on Lithium-ion Battery Prices"
$ws1.Range("B6").Value = "https://about.bnef.com/blog/behind-scenes-take-lithium-ion-battery-prices/"

# Drop the old "Abstract" label and the learning-rate note - no longer applicable
$ws1.Range("B7").ClearContents()
$ws1.Range("A9").ClearContents()

# New footer formatting cell touched while editing the sheet
$ws1.Range("D14").Interior.ColorIndex = -4142

# --- "PDiBCpDoC" sheet: hardcode the new BNEF-sourced learning rate ---
$ws2.Range("B2").Value = 0.18

# --- restore selections / active sheet to match the saved view state ---
[void]$ws2.Range("H30").Select()
[void]$ws1.Range("B11").Select()
